$wb = $excel.ActiveWorkbook

# --- Constants sheet: the previously blank row 2 was deleted, shifting all
#     subsequent rows up by one. ---
$wsConstants = $wb.Worksheets.Item("Constants")
$wsConstants.Rows.Item(2).Delete()
$wsConstants.Activate()
$wsConstants.Range("A12").Select()

# --- Assets sheet: leave its stored selection at A2 (no longer the active tab). ---
$wsAssets = $wb.Worksheets.Item("Assets")
$wsAssets.Activate()
$wsAssets.Range("A2").Select()

# --- Settings sheet: the OrchestratorQueueName cell was retyped (cleared and
#     re-entered), and the sheet became the active tab with the cursor on A6. ---
$wsSettings = $wb.Worksheets.Item("Settings")
$wsSettings.Range("A2").ClearContents()
$wsSettings.Range("A2").Value = "OrchestratorQueueName"
$wsSettings.Activate()
$wsSettings.Range("A6").Select()
